$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 835 ("等々" entry), shifting all subsequent rows up by one.
$ws.Rows.Item(835).Delete()
